$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 80.87054666666667
$ws.Range("H2").Value = 242.61164
$ws.Range("I2").Value = 0.7161501349062054
$ws.Range("J2").Value = 0.7161501349062055
$ws.Range("O2").Value = 0.01124317172993818
$ws.Range("P2").Value = 0.01124317172993818
$ws.Range("Q2").Value = 10.69467153023555
$ws.Range("R2").Value = 96.25204377212
$ws.Range("S2").Value = 0.008051798951168865
$ws.Range("T2").Value = 0.008051798951168867

# Row 3
$ws.Range("G3").Value = 80.87054666666667
$ws.Range("H3").Value = 242.61164
$ws.Range("I3").Value = 0.7161501349062054
$ws.Range("J3").Value = 0.7161501349062055
$ws.Range("M3").Value = 11.62994666666667
$ws.Range("N3").Value = 34.88984
$ws.Range("O3").Value = 0.9887568282700618
$ws.Range("P3").Value = 0.9887568282700618
$ws.Range("Q3").Value = 940.5201446375112
$ws.Range("R3").Value = 8464.681301737599
$ws.Range("S3").Value = 0.7080983359550366
$ws.Range("T3").Value = 0.7080983359550367

# Row 4
$ws.Range("I4").Value = 0.04039551233681073
$ws.Range("J4").Value = 0.04039551233681073
$ws.Range("O4").Value = 0.01124317172993818
$ws.Range("P4").Value = 0.01124317172993818
$ws.Range("S4").Value = 0.0004541736823215994
$ws.Range("T4").Value = 0.0004541736823215996

# Row 5
$ws.Range("I5").Value = 0.04039551233681073
$ws.Range("J5").Value = 0.04039551233681073
$ws.Range("M5").Value = 11.62994666666667
$ws.Range("N5").Value = 34.88984
$ws.Range("O5").Value = 0.9887568282700618
$ws.Range("P5").Value = 0.9887568282700618
$ws.Range("Q5").Value = 53.05143608008889
$ws.Range("R5").Value = 477.4629247208
$ws.Range("S5").Value = 0.03994133865448913
$ws.Range("T5").Value = 0.03994133865448914

# Row 6
$ws.Range("G6").Value = 27.01376833333333
$ws.Range("H6").Value = 81.04130499999999
$ws.Range("I6").Value = 0.2392207624857774
$ws.Range("J6").Value = 0.2392207624857774
$ws.Range("O6").Value = 0.01124317172993818
$ws.Range("P6").Value = 0.01124317172993818
$ws.Range("Q6").Value = 3.572417784062777
$ws.Range("R6").Value = 32.151760056565
$ws.Range("S6").Value = 0.002689600113994349
$ws.Range("T6").Value = 0.002689600113994349

# Row 7
$ws.Range("G7").Value = 27.01376833333333
$ws.Range("H7").Value = 81.04130499999999
$ws.Range("I7").Value = 0.2392207624857774
$ws.Range("J7").Value = 0.2392207624857774
$ws.Range("M7").Value = 11.62994666666667
$ws.Range("N7").Value = 34.88984
$ws.Range("O7").Value = 0.9887568282700618
$ws.Range("P7").Value = 0.9887568282700618
$ws.Range("Q7").Value = 314.1686849823556
$ws.Range("R7").Value = 2827.5181648412
$ws.Range("S7").Value = 0.236531162371783
$ws.Range("T7").Value = 0.2365311623717831

# Row 8
$ws.Range("G8").Value = 0.4780740000000001
$ws.Range("H8").Value = 1.434222
$ws.Range("I8").Value = 0.004233590271206475
$ws.Range("J8").Value = 0.004233590271206476
$ws.Range("O8").Value = 0.01124317172993818
$ws.Range("P8").Value = 0.01124317172993818
$ws.Range("Q8").Value = 0.063222577414
$ws.Range("R8").Value = 0.5690031967260001
$ws.Range("S8").Value = [double]"4.759898245336997E-05"
$ws.Range("T8").Value = [double]"4.759898245336999E-05"

# Row 9
$ws.Range("G9").Value = 0.4780740000000001
$ws.Range("H9").Value = 1.434222
$ws.Range("I9").Value = 0.004233590271206475
$ws.Range("J9").Value = 0.004233590271206476
$ws.Range("M9").Value = 11.62994666666667
$ws.Range("N9").Value = 34.88984
$ws.Range("O9").Value = 0.9887568282700618
$ws.Range("P9").Value = 0.9887568282700618
$ws.Range("Q9").Value = 5.559975122720001
$ws.Range("R9").Value = 50.03977610448
$ws.Range("S9").Value = 0.004185991288753106
$ws.Range("T9").Value = 0.004185991288753107
